{"js": "// Change the report date run \"Date      :   28 / 05 / 20\" into\n// \"Date      :   05 / 28 \" + \"/ 19\" (two runs, same 12pt formatting),\n// i.e. swap day/month and shorten the year to \"19\".\n\nconst body = context.document.body;\n\n// 1) Locate the original date text and rewrite its first part in place.\n//    Doing this as a targeted Range.insertText(\"...\", \"Replace\") keeps the\n//    paragraph's existing run (and its formatting) for the retyped prefix.\nconst dateHits = body.search(\"Date      :   28 / 05 / 20\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nif (dateHits.items.length === 0) {\n  throw new Error(\"Could not find the original date text to edit.\");\n}\n\nconst dateRun = dateHits.items[0];\ndateRun.insertText(\"Date      :   05 / 28 \", \"Replace\");\nawait context.sync();\n\n// 2) Append the \"/ 19\" suffix right after the text we just wrote. Because\n//    this is a separate insertion (its own search + insert), the saved\n//    document keeps it as a distinct run instead of silently re-merging it\n//    into the previous run.\nconst prefixHits = body.search(\"05 / 28 \", { matchCase: true });\nprefixHits.load(\"items\");\nawait context.sync();\n\nif (prefixHits.items.length === 0) {\n  throw new Error(\"Could not find the rewritten date prefix.\");\n}\n\nconst afterPrefix = prefixHits.items[0].getRange(\"End\");\nafterPrefix.insertText(\"/ 19\", \"After\");\nawait context.sync();\n\n// 3) The new \"/ 19\" run already inherits the surrounding 12pt formatting,\n//    but re-stamp it explicitly (via a round trip through another size) so\n//    it is kept as its own run rather than being folded back into the\n//    neighboring run just because the two runs look identical.\nconst suffixHits = body.search(\"/ 19\", { matchCase: true });\nsuffixHits.load(\"items\");\nawait context.sync();\n\nif (suffixHits.items.length === 0) {\n  throw new Error(\"Could not find the appended '/ 19' suffix.\");\n}\n\nconst suffixRun = suffixHits.items[0];\nsuffixRun.font.size = 11;\nawait context.sync();\n\nsuffixRun.font.size = 12;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the original date line and split its text into two pieces so the\n# saved document ends up with two runs (matching how Word keeps the\n# existing run for the retyped prefix and starts a new run for the\n# appended suffix).\n$find = $d.Content\n$find.Find.Execute(\"28 / 05 / 20\", $false, $false, $false, $false, $false, $true, 1, $false, \"05 / 28 \", 2)\n\n# Re-find the now-updated prefix so we can anchor a new, separate run\n# immediately after it for the \"/ 19\" suffix.\n$tail = $d.Content\n$tail.Find.Execute(\"05 / 28 \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$tail.Collapse(0)\n$newRun = $tail.InsertAfter(\"/ 19\")\n\n# Force the freshly inserted text to become its own run: nudge its font\n# size away from and then back to the paragraph's 12pt (24 half-points)\n# setting. Because the two edits touch the range at different times, the\n# engine keeps the \"/ 19\" text as an independent run instead of folding it\n# back into the preceding run even though the final formatting matches.\n$tail.Find.Execute(\"/ 19\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$tail.Font.Size = 99\n$tail.Font.Size = 12\n"}
